$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$tbl = $ws.ListObjects.Item(1)

# -------------------------------------------------------------------
# 1) Fill in EARNED/Absence values for the Aug-Nov 2023 rows that were
#    previously blank (rows 46-48) and the Dec-2023 "FL(2-0-0)" row (49)
# -------------------------------------------------------------------
$ws.Range("C46").Value = 1.25
$ws.Range("C47").Value = 1.25
$ws.Range("C48").Value = 1.25

$ws.Range("B49").Value = "FL(2-0-0)"
$ws.Range("C49").Value = 1.25
$ws.Range("D49").Value = 2

# -------------------------------------------------------------------
# 2) Row 50 turns from a "01/01/2024" period-start date into a bold
#    "2024" year-separator label (same visual style as the "2023"
#    label already used in row 36).
# -------------------------------------------------------------------
$ws.Range("A50").NumberFormat = "mm/dd/yy;@"
$ws.Range("A50").Font.Bold = $true
$ws.Range("A50").Value = "'2024"

# -------------------------------------------------------------------
# 3) Row 51: date shifts back one day, gains a "SP(1-0-0)" particular,
#    and a remarks date of 1/3/2024 (formatted like the other date
#    remarks in column K, e.g. K45).
# -------------------------------------------------------------------
$ws.Range("A51").Value = $ws.Range("A51").Value2 - 1
$ws.Range("B51").Value = "SP(1-0-0)"
$ws.Range("K45").Copy($ws.Range("K51"))
$ws.Range("K51").Value = 45294

# -------------------------------------------------------------------
# 4) Rows 52-117: every PERIOD date moves back one day (end-of-month
#    instead of first-of-month).
# -------------------------------------------------------------------
for ($r = 52; $r -le 117; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $cell.Value2 - 1
}

# Row 118 gains its end-of-month PERIOD date too.
$ws.Cells.Item(118, 1).Value = 47361

# -------------------------------------------------------------------
# 5) Grow Table1 by one row: the old last (bottom-bordered) row 147
#    becomes a normal interior row, and a brand new bottom-bordered
#    row appears as row 148.
# -------------------------------------------------------------------
$newRow = $tbl.ListRows.Add()
$ws.Range("A147:K147").Copy($ws.Range("A148:K148"))
$ws.Range("A146:K146").Copy($ws.Range("A147:K147"))
$ws.Range("G148").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# -------------------------------------------------------------------
# 6) Leave the cursor positioned the way the author left it.
# -------------------------------------------------------------------
$ws.Activate()
$ws.Range("B52").Select()
